$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# scrapedAt (H) / lastSeenAt (I) timestamps refreshed by the latest scrape run
$timestamps = @(
    @(2, "2026-02-03T14:45:23.454Z", "2026-02-03T14:45:23.493Z"),
    @(3, "2026-02-03T14:45:23.455Z", "2026-02-03T14:45:23.511Z"),
    @(4, "2026-02-03T14:45:23.455Z", "2026-02-03T14:45:23.514Z"),
    @(5, "2026-02-03T14:45:23.455Z", "2026-02-03T14:45:23.517Z"),
    @(6, "2026-02-03T14:45:23.455Z", "2026-02-03T14:45:23.519Z"),
    @(7, "2026-02-03T14:45:23.455Z", "2026-02-03T14:45:23.523Z"),
    @(8, "2026-02-03T14:45:23.455Z", "2026-02-03T14:45:23.524Z"),
    @(9, "2026-02-03T14:45:23.456Z", "2026-02-03T14:45:23.528Z"),
    @(10, "2026-02-03T14:45:23.456Z", "2026-02-03T14:45:23.530Z"),
    @(11, "2026-02-03T14:45:23.456Z", "2026-02-03T14:45:23.533Z"),
    @(12, "2026-02-03T14:45:23.456Z", "2026-02-03T14:45:23.535Z"),
    @(13, "2026-02-03T14:45:23.456Z", "2026-02-03T14:45:23.538Z"),
    @(14, "2026-02-03T14:45:23.456Z", "2026-02-03T14:45:23.543Z"),
    @(15, "2026-02-03T14:45:23.456Z", "2026-02-03T14:45:23.545Z"),
    @(16, "2026-02-03T14:45:23.457Z", "2026-02-03T14:45:23.547Z"),
    @(17, "2026-02-03T14:45:23.457Z", "2026-02-03T14:45:23.549Z"),
    @(18, "2026-02-03T14:45:23.457Z", "2026-02-03T14:45:23.551Z"),
    @(19, "2026-02-03T14:45:23.457Z", "2026-02-03T14:45:23.553Z"),
    @(20, "2026-02-03T14:45:23.457Z", "2026-02-03T14:45:23.555Z"),
    @(21, "2026-02-03T14:45:23.457Z", "2026-02-03T14:45:23.557Z"),
    @(22, "2026-02-03T14:45:23.457Z", "2026-02-03T14:45:23.559Z"),
    @(23, "2026-02-03T14:45:23.458Z", "2026-02-03T14:45:23.561Z"),
    @(24, "2026-02-03T14:45:23.458Z", "2026-02-03T14:45:23.564Z"),
    @(25, "2026-02-03T14:45:23.458Z", "2026-02-03T14:45:23.566Z"),
    @(26, "2026-02-03T14:45:23.458Z", "2026-02-03T14:45:23.569Z"),
    @(27, "2026-02-03T14:45:23.459Z", "2026-02-03T14:45:23.575Z"),
    @(28, "2026-02-03T14:45:23.459Z", "2026-02-03T14:45:23.578Z"),
    @(29, "2026-02-03T14:45:23.459Z", "2026-02-03T14:45:23.582Z"),
    @(30, "2026-02-03T14:45:23.459Z", "2026-02-03T14:45:23.585Z"),
    @(31, "2026-02-03T14:45:23.460Z", "2026-02-03T14:45:23.588Z"),
    @(32, "2026-02-03T14:45:23.460Z", "2026-02-03T14:45:23.591Z"),
    @(33, "2026-02-03T14:45:23.460Z", "2026-02-03T14:45:23.601Z"),
    @(34, "2026-02-03T14:45:23.460Z", "2026-02-03T14:45:23.604Z"),
    @(35, "2026-02-03T14:45:23.460Z", "2026-02-03T14:45:23.607Z"),
    @(36, "2026-02-03T14:45:23.460Z", "2026-02-03T14:45:23.610Z"),
    @(37, "2026-02-03T14:45:23.460Z", "2026-02-03T14:45:23.612Z"),
    @(38, "2026-02-03T14:45:23.460Z", "2026-02-03T14:45:23.615Z"),
    @(39, "2026-02-03T14:45:23.460Z", "2026-02-03T14:45:23.618Z"),
    @(40, "2026-02-03T14:45:23.460Z", "2026-02-03T14:45:23.621Z"),
    @(41, "2026-02-03T14:45:23.460Z", "2026-02-03T14:45:23.624Z"),
    @(42, "2026-02-03T14:45:23.461Z", "2026-02-03T14:45:23.626Z"),
    @(43, "2026-02-03T14:45:23.461Z", "2026-02-03T14:45:23.628Z"),
    @(44, "2026-02-03T14:45:23.461Z", "2026-02-03T14:45:23.630Z"),
    @(45, "2026-02-03T14:45:23.461Z", "2026-02-03T14:45:23.633Z"),
    @(46, "2026-02-03T14:45:23.461Z", "2026-02-03T14:45:23.636Z")
)

foreach ($row in $timestamps) {
    $r = $row[0]
    $ws.Cells.Item($r, 8).Value = $row[1]
    $ws.Cells.Item($r, 9).Value = $row[2]
}
